# Update cryptos list values (price & 1h volume change) per latest data pull
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.806.72"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "'3.259.61"
$ws.Range("E3").Value = "  +2.68%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'604.71"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "'158.23"
$ws.Range("E6").Value = "  +2.19%  "

$ws.Range("D8").Value = "'3.260.46"
$ws.Range("E8").Value = "  +2.75%  "

$ws.Range("D9").Value = "'0.549"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +3.10%  "

$ws.Range("D11").Value = "'5.88"
$ws.Range("E11").Value = "  +3.76%  "

$ws.Range("D12").Value = "'0.507"
$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("E13").Value = "  +3.09%  "

$ws.Range("D14").Value = "'39.44"
$ws.Range("E14").Value = "  +2.13%  "

$ws.Range("D15").Value = "'3.799.37"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").Value = "'66.856.06"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").Value = "'7.39"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").Value = "'3.268.42"
$ws.Range("E18").Value = "  +2.84%  "

$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").Value = "'508.83"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").Value = "'15.44"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").Value = "'0.755"
$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("D23").Value = "'8.10"
$ws.Range("E23").Value = "  -0.50%  "

$ws.Range("D24").Value = "'14.79"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").Value = "'86.33"
$ws.Range("E25").Value = "  +2.30%  "

$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +87.49%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").Value = "'3.02"
$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("D29").Value = "'9.11"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'6.93"
$ws.Range("E31").Value = "  -1.18%  "

$ws.Range("D32").Value = "'2.88"
$ws.Range("E32").Value = "  -6.41%  "

$ws.Range("D33").Value = "'28.33"
$ws.Range("E33").Value = "  +1.14%  "

$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").Value = "'1.16"
$ws.Range("E35").Value = "  -3.69%  "

$ws.Range("D36").Value = "'6.44"
$ws.Range("E36").Value = "  -1.27%  "

$ws.Range("D37").Value = "'0.0₃0807"
$ws.Range("E37").Value = "  +20.09%  "

$ws.Range("D38").Value = "'3.36"
$ws.Range("E38").Value = "  +19.92%  "

$ws.Range("D39").Value = "'55.63"
$ws.Range("E39").Value = "  +1.75%  "

$ws.Range("D40").Value = "'496.35"
$ws.Range("E40").Value = "  -2.23%  "

$ws.Range("D41").Value = "'0.0430"
$ws.Range("E41").Value = "  +2.11%  "

$ws.Range("D42").Value = "'0.128"
$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("D43").Value = "'8.80"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").Value = "'0.297"
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").Value = "'2.49"
$ws.Range("E45").Value = "  +3.12%  "

$ws.Range("D46").Value = "'2.947.80"
$ws.Range("E46").Value = "  +3.78%  "

$ws.Range("D47").Value = "'28.65"
$ws.Range("E47").Value = "  +1.54%  "

$ws.Range("E48").Value = "  +1.74%  "

$ws.Range("D49").Value = "'0.119"
$ws.Range("E49").Value = "  +2.21%  "

$ws.Range("D51").Value = "'2.55"
$ws.Range("E51").Value = "  -0.12%  "
